$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unlock it so the cell values below can be updated,
# then restore protection once the edits are in place.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A11)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."
# Re-fit row 11 so the embedded line break doesn't leave a stray explicit row height
$ws.Rows.Item(11).AutoFit()

# Update the Weight (D) and Percent Change (E) figures for each holding row
$ws.Range("D2").Value = 0.4939161058383833
$ws.Range("E2").Value = -0.0006968641114981633

$ws.Range("D3").Value = 0.2476925908086156
$ws.Range("E3").Value = -0.001039038147543581

$ws.Range("D4").Value = 0.09867507434849365
$ws.Range("E4").Value = 0.001373454863278889

$ws.Range("D5").Value = 0.1015856514007556
$ws.Range("E5").Value = -0.0006737895851380848

$ws.Range("D6").Value = 0.02994767297575263
$ws.Range("E6").Value = -0.0005877167205408096

$ws.Range("D7").Value = 0.02818290462799939
$ws.Range("E7").Value = -0.0006869704602702198

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = -0.000571437623241744

# Restore the original sheet protection state
$ws.Protect()
